# Correct misleading values in mapping schemes; revised area and cost
# assumptions for all occupancies; revised count assumptions for
# non-residential, on the "Dwellings_buildings" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dwellings_buildings")

# Classification label renames (col B)
$ws.Range("B3").Value = "Professional and technical services"
$ws.Range("B4").Value = "All other services"
$ws.Range("B6").Value = "Professional and technical services"
$ws.Range("B7").Value = "All other services"
$ws.Range("B9").Value = "Professional and technical services"
$ws.Range("B10").Value = "All other services"
$ws.Range("B12").Value = "Professional and technical services"
$ws.Range("B13").Value = "All other services"
$ws.Range("B14").Value = "Professional and technical services"
$ws.Range("B15").Value = "Professional and technical services"
$ws.Range("B16").Value = "All other services"
$ws.Range("B17").Value = "Professional and technical services"
$ws.Range("B18").Value = "Professional and technical services"
$ws.Range("B19").Value = "All other services"

# average_area (col C)
$ws.Range("C2").Value = 130
$ws.Range("C3").Value = 130
$ws.Range("C4").Value = 130
$ws.Range("C5").Value = 260
$ws.Range("C6").Value = 260
$ws.Range("C7").Value = 260
$ws.Range("C8").Value = 450
$ws.Range("C9").Value = 450
$ws.Range("C10").Value = 450
$ws.Range("C11").Value = 900
$ws.Range("C12").Value = 900
$ws.Range("C13").Value = 900
$ws.Range("C14").Value = 1200
$ws.Range("C15").Value = 1200
$ws.Range("C16").Value = 1200
$ws.Range("C17").Value = 3200
$ws.Range("C18").Value = 3200
$ws.Range("C19").Value = 3200

# average_unit_cost (col D)
$ws.Range("D2").Value = 453.2997557003257
$ws.Range("D3").Value = 520.6709120521173
$ws.Range("D4").Value = 469.9346091205211
$ws.Range("D5").Value = 453.2997557003257
$ws.Range("D6").Value = 520.6709120521173
$ws.Range("D7").Value = 469.9346091205211
$ws.Range("D8").Value = 453.2997557003257
$ws.Range("D9").Value = 520.6709120521173
$ws.Range("D10").Value = 469.9346091205211
$ws.Range("D11").Value = 453.2997557003257
$ws.Range("D12").Value = 520.6709120521173
$ws.Range("D13").Value = 469.9346091205211
$ws.Range("D14").Value = 682.028990228013
$ws.Range("D15").Value = 453.2997557003257
$ws.Range("D16").Value = 469.9346091205211
$ws.Range("D17").Value = 682.028990228013
$ws.Range("D18").Value = 453.2997557003257
$ws.Range("D19").Value = 469.9346091205211

# dwellings_per_building (col F)
$ws.Range("F5").Value = 2
$ws.Range("F6").Value = 2
$ws.Range("F7").Value = 2
$ws.Range("F8").Value = 3
$ws.Range("F9").Value = 3
$ws.Range("F10").Value = 3
$ws.Range("F11").Value = 5
$ws.Range("F12").Value = 5
$ws.Range("F13").Value = 5
$ws.Range("F14").Value = 5
$ws.Range("F15").Value = 5
$ws.Range("F16").Value = 5
$ws.Range("F17").Value = 10
$ws.Range("F18").Value = 10
$ws.Range("F19").Value = 10
